$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 111934050
$ws.Range("B5").Value = 93289
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2170
$ws.Range("F5").Value = "Flagellkvastmossa"
$ws.Range("G5").Value = "Dicranum flagellare"
$ws.Range("H5").Value = "Hedw."
$ws.Range("K5").Value = "med groddkorn"
$ws.Range("P5").Value = "Skogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q5").Value = 413637.9321653559
$ws.Range("R5").Value = 6587076.603947581
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Värmland"
$ws.Range("U5").Value = "Karlstad"
$ws.Range("V5").Value = "Värmland"
$ws.Range("W5").Value = "Grava"
$ws.Range("Y5").Value = "'2023-09-06"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").Value = "'2023-09-06"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "Carl A Andersson"
$ws.Range("AX5").Value = "Carl A Andersson"

$ws.Range("A6").Value = 111934086
$ws.Range("B6").Value = 90689
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 5966
$ws.Range("F6").Value = "Motaggsvamp"
$ws.Range("G6").Value = "Sarcodon squamosus"
$ws.Range("H6").Value = "(Schaeff.) Quél."
$ws.Range("I6").Value = "'1"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("P6").Value = "Tallskogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q6").Value = 413681.2082122188
$ws.Range("R6").Value = 6586805.223123537
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Värmland"
$ws.Range("U6").Value = "Karlstad"
$ws.Range("V6").Value = "Värmland"
$ws.Range("W6").Value = "Grava"
$ws.Range("Y6").Value = "'2023-09-06"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").Value = "'2023-09-06"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Carl A Andersson"
$ws.Range("AX6").Value = "Carl A Andersson"

$ws.Range("A7").Value = 111934066
$ws.Range("B7").Value = 93289
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 2170
$ws.Range("F7").Value = "Flagellkvastmossa"
$ws.Range("G7").Value = "Dicranum flagellare"
$ws.Range("H7").Value = "Hedw."
$ws.Range("K7").Value = "med groddkorn"
$ws.Range("P7").Value = "Skogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q7").Value = 413590.3038565172
$ws.Range("R7").Value = 6586912.201658082
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Värmland"
$ws.Range("U7").Value = "Karlstad"
$ws.Range("V7").Value = "Värmland"
$ws.Range("W7").Value = "Grava"
$ws.Range("Y7").Value = "'2023-09-06"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").Value = "'2023-09-06"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = "Carl A Andersson"
$ws.Range("AX7").Value = "Carl A Andersson"

$ws.Range("A8").Value = 111934059
$ws.Range("B8").Value = 93289
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 2170
$ws.Range("F8").Value = "Flagellkvastmossa"
$ws.Range("G8").Value = "Dicranum flagellare"
$ws.Range("H8").Value = "Hedw."
$ws.Range("K8").Value = "med groddkorn"
$ws.Range("P8").Value = "Skogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q8").Value = 413639.6308819132
$ws.Range("R8").Value = 6586793.951973591
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Värmland"
$ws.Range("U8").Value = "Karlstad"
$ws.Range("V8").Value = "Värmland"
$ws.Range("W8").Value = "Grava"
$ws.Range("Y8").Value = "'2023-09-06"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").Value = "'2023-09-06"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AC8").Value = "Rätt riklig längs stigen"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = "Carl A Andersson"
$ws.Range("AX8").Value = "Carl A Andersson"

Write-Host "done"
